$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the value must be forced to
# Text so Excel does not silently reinterpret a numeric-looking string (e.g.
# "239.26" or "5.250") as a Number and drop significant trailing zeros.
$updates = @(
    @{ Cell = 'D2'; Value = '25.998.97'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +0.66%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.739.87'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +0.09%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '239.26'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +3.18%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '1.002'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +0.16%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.5283'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  +2.15%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.2727'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -2.72%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.06159'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +0.74%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '1.743.80'; ForceText = $false }
    @{ Cell = 'E10'; Value = '  -0.38%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.07181'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +2.08%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '15.05'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -1.62%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '0.6400'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  -0.88%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '4.604'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  +1.73%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '77.50'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +0.71%  '; ForceText = $false }
    @{ Cell = 'E16'; Value = '  +0.09%  '; ForceText = $false }
    @{ Cell = 'E17'; Value = '  +0.09%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '26.009.44'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +0.73%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '11.76'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +2.34%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.000006748'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +2.23%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '1.963.54'; ForceText = $false }
    @{ Cell = 'E21'; Value = '  -0.51%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '4.349'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  +5.05%  '; ForceText = $false }
    @{ Cell = 'E23'; Value = '  -0.64%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '5.250'; ForceText = $true }
    @{ Cell = 'D25'; Value = '140.16'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +0.47%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '1.502'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -0.67%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '15.25'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +1.25%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '1.763'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -2.66%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '105.76'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  +3.57%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '0.08389'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +0.67%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '3.829'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +4.03%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '3.642'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +6.33%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '0.04589'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +2.12%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '2.656'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +1.84%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '0.9920'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +0.75%  '; ForceText = $false }
    @{ Cell = 'E36'; Value = '  +1.38%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '2.695'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  +1.63%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.01602'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +1.19%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '1.926'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -0.81%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '1.002'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +0.22%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '98.70'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -1.96%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.3880'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +0.93%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.7447'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  +2.10%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '4.943'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -0.67%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.1142'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +1.97%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '0.05310'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -1.83%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '6.202'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -1.23%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '54.64'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +2.91%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '30.71'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  +2.43%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '7.515'; ForceText = $true }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Cell = 'E50'; Value = '  -2.24%  '; ForceText = $false }
    @{ Cell = 'B50'; Value = 'EnergySwap'; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.3442'; ForceText = $true }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; ForceText = $false }
    @{ Cell = 'E51'; Value = '  +1.38%  '; ForceText = $false }
    @{ Cell = 'B51'; Value = 'Decentraland'; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Temporarily mark the cell as Text so the numeric-looking string is
        # stored verbatim, then clear the formatting again so the cell keeps
        # its original (default) style, matching the source edit exactly.
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.ClearFormats()
    } else {
        $range.Value = $u.Value
    }
}
